$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.628.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'2.089.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'342.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").Value = "'1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "'0.5162"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.97%  "
$ws.Range("D8").Value = "'0.4391"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.57%  "
$ws.Range("D9").Value = "'0.09256"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("D10").Value = "'51.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.26%  "
$ws.Range("D11").Value = "'1.177"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "'25.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.72%  "
$ws.Range("D13").Value = "'2.089.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "'6.742"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").Value = "'8.179"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").Value = "'100.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "'0.00001157"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "'1.009"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "'21.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.29%  "
$ws.Range("D20").Value = "'0.06628"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").Value = "'1.008"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").Value = "'6.182"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").Value = "'29.675.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.40%  "
$ws.Range("D24").Value = "'12.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("D25").Value = "'2.312"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.89%  "
$ws.Range("D26").Value = "'2.334.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("D27").Value = "'21.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("D28").Value = "'163.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("D29").Value = "'2.519"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").Value = "'132.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("D31").Value = "'1.137"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.17%  "
$ws.Range("E32").Value = "  -2.28%  "
$ws.Range("D33").Value = "'1.630"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("D34").Value = "'6.185"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("D35").Value = "'3.958"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("D36").Value = "'6.030"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.38%  "
$ws.Range("D37").Value = "'10.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("D38").Value = "'0.02569"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.72%  "
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("D40").Value = "'12.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").Value = "'0.2240"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("D42").Value = "'0.6824"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").Value = "'1.291"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").Value = "'0.6600"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("D45").Value = "'14.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.06%  "
$ws.Range("D46").Value = "'2.313"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").Value = "'3.604"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.52%  "
$ws.Range("D48").Value = "'1.217"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("D49").Value = "'0.00000000338"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.43%  "
$ws.Range("D50").Value = "'81.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").Value = "'1.166"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.22%  "

Write-Output "Applied cryptos update"